$d = $word.ActiveDocument

# 1. Shorten the "Curso (semestre ideal)" line: drop EF (9), EM (8), EB (8)
$d.Content.Find.Execute(
    "Curso (semestre ideal): EF (9), EM (8), EB (8), EP (10), EQD (8), EQN (11)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EP (10), EQD (8), EQN (11)", 2)

# 2. Remove the trailing "Requisitos" section (heading paragraph + the
#    "LOB1008 - ..." paragraph that follows it), which were the last two
#    paragraphs of the document body.
$count = $d.Paragraphs.Count
$d.Paragraphs.Item($count).Range.Delete()
$d.Paragraphs.Item($count - 1).Range.Delete()
